# Auto-generated edit script applying Famfrit_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 463.33334
$ws.Range("I18").Value = 463.33334
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 463.33334
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -179.33334
$ws.Range("H40").Value = 1434327.8
$ws.Range("I40").Value = 1672049
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 1672049
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -1671874
$ws.Range("N40").Value = -8350
$ws.Range("H74").Value = 5779.467
$ws.Range("I74").Value = 5223
$ws.Range("J74").Value = 5981.8184
$ws.Range("K74").Value = 5223
$ws.Range("L74").Value = 5981.8184
$ws.Range("M74").Value = -4287
$ws.Range("N74").Value = -7853.8184
$ws.Range("H77").Value = 5779.467
$ws.Range("I77").Value = 5223
$ws.Range("J77").Value = 5981.8184
$ws.Range("K77").Value = 26115
$ws.Range("L77").Value = 29909.092
$ws.Range("M77").Value = -21435
$ws.Range("N77").Value = -39269.092
$ws.Range("H98").Value = 3075.625
$ws.Range("I98").Value = 2767.5
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 2767.5
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = -1269.5
$ws.Range("N98").Value = -6996
$ws.Range("H113").Value = 10997.6
$ws.Range("I113").Value = 9997
$ws.Range("J113").Value = 15000
$ws.Range("K113").Value = 9997
$ws.Range("L113").Value = 15000
$ws.Range("M113").Value = -6743
$ws.Range("N113").Value = -21508
$ws.Range("H122").Value = 3075.625
$ws.Range("I122").Value = 2767.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8302.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5852.5
$ws.Range("N122").Value = -16900
$ws.Range("H137").Value = 4342.9414
$ws.Range("I137").Value = 4708.8184
$ws.Range("J137").Value = 3672.1667
$ws.Range("K137").Value = 14126.4552
$ws.Range("L137").Value = 11016.5001
$ws.Range("M137").Value = -11576.4552
$ws.Range("N137").Value = -16116.5001
$ws.Range("H138").Value = 3990.6099
$ws.Range("I138").Value = 1114.2693
$ws.Range("J138").Value = 8976.267
$ws.Range("K138").Value = 3342.8079
$ws.Range("L138").Value = 26928.801
$ws.Range("M138").Value = 1797.1921
$ws.Range("N138").Value = -37208.801
$ws.Range("H141").Value = 1258.9773
$ws.Range("I141").Value = 1186.5853
$ws.Range("J141").Value = 2248.3333
$ws.Range("K141").Value = 3559.7559
$ws.Range("L141").Value = 6744.999899999999
$ws.Range("M141").Value = 1620.2441
$ws.Range("N141").Value = -17104.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5907.1646
$ws.Range("I32").Value = 4467.691
$ws.Range("J32").Value = 10305.556
$ws.Range("K32").Value = 4467.691
$ws.Range("L32").Value = 10305.556
$ws.Range("M32").Value = -4180.691
$ws.Range("N32").Value = -10879.556
$ws.Range("H132").Value = 37136268
$ws.Range("I132").Value = 14223.723
$ws.Range("J132").Value = 111380350
$ws.Range("K132").Value = 42671.169
$ws.Range("L132").Value = 334141050
$ws.Range("M132").Value = -40141.169
$ws.Range("N132").Value = -334146110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H48").Value = 269999
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 269999
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 269999
$ws.Range("N48").Value = -270829

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1375.875
$ws.Range("I16").Value = 1102.2
$ws.Range("J16").Value = 1832
$ws.Range("K16").Value = 1102.2
$ws.Range("L16").Value = 1832
$ws.Range("M16").Value = -815.2
$ws.Range("N16").Value = -2406
$ws.Range("H29").Value = 3500
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 5000
$ws.Range("M29").Value = -1707
$ws.Range("N29").Value = -5586
$ws.Range("H31").Value = 4842.1304
$ws.Range("I31").Value = 3088.2942
$ws.Range("J31").Value = 9811.333000000001
$ws.Range("K31").Value = 3088.2942
$ws.Range("L31").Value = 9811.333000000001
$ws.Range("M31").Value = -2793.2942
$ws.Range("N31").Value = -10401.333
$ws.Range("H34").Value = 4842.1304
$ws.Range("I34").Value = 3088.2942
$ws.Range("J34").Value = 9811.333000000001
$ws.Range("K34").Value = 3088.2942
$ws.Range("L34").Value = 9811.333000000001
$ws.Range("M34").Value = -2886.2942
$ws.Range("N34").Value = -10215.333
$ws.Range("H58").Value = 2126.5454
$ws.Range("I58").Value = 707.3077
$ws.Range("J58").Value = 4176.5557
$ws.Range("K58").Value = 707.3077
$ws.Range("L58").Value = 4176.5557
$ws.Range("M58").Value = -504.3077
$ws.Range("N58").Value = -4582.5557
$ws.Range("H113").Value = 1375.875
$ws.Range("I113").Value = 1102.2
$ws.Range("J113").Value = 1832
$ws.Range("K113").Value = 1102.2
$ws.Range("L113").Value = 1832
$ws.Range("M113").Value = 1067.8
$ws.Range("N113").Value = -6172
$ws.Range("H122").Value = 1058.2
$ws.Range("I122").Value = 1058.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3174.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -724.6000000000004
$ws.Range("H132").Value = 9807.429
$ws.Range("I132").Value = 9775.333000000001
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 29325.999
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -26795.999
$ws.Range("N132").Value = -35060
$ws.Range("H136").Value = 2126.5454
$ws.Range("I136").Value = 707.3077
$ws.Range("J136").Value = 4176.5557
$ws.Range("K136").Value = 2121.9231
$ws.Range("L136").Value = 12529.6671
$ws.Range("M136").Value = 428.0769
$ws.Range("N136").Value = -17629.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42928080
$ws.Range("I4").Value = 56316590
$ws.Range("J4").Value = 25714286
$ws.Range("K4").Value = 168949770
$ws.Range("L4").Value = 77142858
$ws.Range("M4").Value = -168949658
$ws.Range("N4").Value = -77143082
$ws.Range("H12").Value = 794
$ws.Range("I12").Value = 228.5
$ws.Range("J12").Value = 1642.25
$ws.Range("K12").Value = 685.5
$ws.Range("L12").Value = 4926.75
$ws.Range("M12").Value = -512.5
$ws.Range("N12").Value = -5272.75
$ws.Range("H56").Value = 7762.2856
$ws.Range("I56").Value = 7762.2856
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 7762.2856
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -7232.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4362.375
$ws.Range("I80").Value = 3724.75
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 3724.75
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -2726.75
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 4362.375
$ws.Range("I83").Value = 3724.75
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 18623.75
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -13631.75
$ws.Range("N83").Value = -34984
$ws.Range("H113").Value = 3799.5945
$ws.Range("I113").Value = 3109.2
$ws.Range("J113").Value = 4611.8237
$ws.Range("K113").Value = 3109.2
$ws.Range("L113").Value = 4611.8237
$ws.Range("M113").Value = -939.1999999999998
$ws.Range("N113").Value = -8951.823700000001
$ws.Range("H122").Value = 2978.5386
$ws.Range("I122").Value = 2717.1
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 8151.299999999999
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -5701.299999999999
$ws.Range("N122").Value = -16450
$ws.Range("H132").Value = 4403.8335
$ws.Range("I132").Value = 3806.5625
$ws.Range("J132").Value = 9182
$ws.Range("K132").Value = 11419.6875
$ws.Range("L132").Value = 27546
$ws.Range("M132").Value = -8889.6875
$ws.Range("N132").Value = -32606

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3411.1
$ws.Range("I7").Value = 2806.3333
$ws.Range("J7").Value = 4822.222
$ws.Range("K7").Value = 2806.3333
$ws.Range("L7").Value = 4822.222
$ws.Range("M7").Value = -2694.3333
$ws.Range("N7").Value = -5046.222
$ws.Range("H40").Value = 3138.1052
$ws.Range("I40").Value = 3072.1667
$ws.Range("J40").Value = 4325
$ws.Range("K40").Value = 3072.1667
$ws.Range("L40").Value = 4325
$ws.Range("M40").Value = -2936.1667
$ws.Range("N40").Value = -4597
$ws.Range("H122").Value = 3558.5278
$ws.Range("I122").Value = 2625.1904
$ws.Range("J122").Value = 4865.2
$ws.Range("K122").Value = 7875.5712
$ws.Range("L122").Value = 14595.6
$ws.Range("M122").Value = -5425.5712
$ws.Range("N122").Value = -19495.6
$ws.Range("H126").Value = 3411.1
$ws.Range("I126").Value = 2806.3333
$ws.Range("J126").Value = 4822.222
$ws.Range("K126").Value = 8418.999899999999
$ws.Range("L126").Value = 14466.666
$ws.Range("M126").Value = -5948.999899999999
$ws.Range("N126").Value = -19406.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 51607.45
$ws.Range("I122").Value = 57158.5
$ws.Range("J122").Value = 1648
$ws.Range("K122").Value = 171475.5
$ws.Range("L122").Value = 4944
$ws.Range("M122").Value = -169025.5
$ws.Range("N122").Value = -9844
